# Commit: Fri, May 29, 2020  2:06:12 PM
#
# 1) The table on slide 16 (the cash-flow summary table) is switched from
#    the deck's one custom table style to a built-in PowerPoint table
#    style ("Medium Style 2 - Accent 1", GUID F550F852-E5EB-457B-B9DB-91D2BB3F5B42).
#
# 2) The presentation's theme (color scheme) is swapped back to the
#    default Office palette, matching the colors that used to live in
#    the deck's secondary ("Office Theme") theme part.

$p = $ppt.ActivePresentation

# --- 1) Re-style the table on slide 16 -------------------------------------
$slide = $p.Slides.Item(16)
for ($i = 1; $i -le $slide.Shapes.Count; $i++) {
    $shape = $slide.Shapes.Item($i)
    if ($shape.HasTable) {
        $shape.Table.ApplyStyle("{F550F852-E5EB-457B-B9DB-91D2BB3F5B42}")
    }
}

# --- 2) Restore the default "Office Theme" color scheme --------------------
$tcs = $p.SlideMaster.Theme.ThemeColorScheme

$tcs.Item(1).RGB  = 0          # dk1      000000
$tcs.Item(2).RGB  = 16777215   # lt1      FFFFFF
$tcs.Item(3).RGB  = 6968388    # dk2      44546A
$tcs.Item(4).RGB  = 15132391   # lt2      E7E6E6
$tcs.Item(5).RGB  = 13998939   # accent1  5B9BD5
$tcs.Item(6).RGB  = 3243501    # accent2  ED7D31
$tcs.Item(7).RGB  = 10855845   # accent3  A5A5A5
$tcs.Item(8).RGB  = 49407      # accent4  FFC000
$tcs.Item(9).RGB  = 12874308   # accent5  4472C4
$tcs.Item(10).RGB = 4697456    # accent6  70AD47
$tcs.Item(11).RGB = 12673797   # hlink    0563C1
$tcs.Item(12).RGB = 7491477    # folHlink 954F72
